$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Text number format on percentage cells first, to prevent Excel from
# auto-converting "NN%" text into a numeric percentage value on assignment.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H45").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-16 23:18:31"
$ws.Range("N2").Value = "0.0 °C 22:44 TU"
$ws.Range("O2").Value = "1.9 °C"
$ws.Range("E3").Value = "2026-02-16 23:18:34"
$ws.Range("E4").Value = "2026-02-16 23:18:36"
$ws.Range("H4").Value = "64%"
$ws.Range("J4").Value = "1012.5 hPa"
$ws.Range("O4").Value = "13.0 °C"
$ws.Range("E5").Value = "2026-02-16 23:18:39"
$ws.Range("H5").Value = "96%"
$ws.Range("N5").Value = "-5.3 °C 22:51 TU"
$ws.Range("O5").Value = "-1.1 °C"
$ws.Range("E6").Value = "2026-02-16 23:18:41"
$ws.Range("K6").Value = "13.1 MJ/m2"
$ws.Range("E7").Value = "2026-02-16 23:18:44"
$ws.Range("H7").Value = "52%"
$ws.Range("E8").Value = "2026-02-16 23:18:46"
$ws.Range("E9").Value = "2026-02-16 23:18:48"
$ws.Range("E10").Value = "2026-02-16 23:18:51"
$ws.Range("O10").Value = "10.5 °C"
$ws.Range("E11").Value = "2026-02-16 23:18:53"
$ws.Range("E12").Value = "2026-02-16 23:18:56"
$ws.Range("H12").Value = "76%"
$ws.Range("O12").Value = "11.3 °C"
$ws.Range("E13").Value = "2026-02-16 23:18:58"
$ws.Range("H13").Value = "75%"
$ws.Range("E14").Value = "2026-02-16 23:19:01"
$ws.Range("H14").Value = "57%"
$ws.Range("E15").Value = "2026-02-16 23:19:03"
$ws.Range("H15").Value = "64%"
$ws.Range("E16").Value = "2026-02-16 23:19:05"
$ws.Range("H16").Value = "79%"
$ws.Range("N16").Value = "-5.3 °C 22:49 TU"
$ws.Range("O16").Value = "-0.5 °C"
$ws.Range("E17").Value = "2026-02-16 23:19:08"
$ws.Range("O17").Value = "5.9 °C"
$ws.Range("E18").Value = "2026-02-16 23:19:10"
$ws.Range("O18").Value = "10.5 °C"
$ws.Range("E19").Value = "2026-02-16 23:19:12"
$ws.Range("H19").Value = "81%"
$ws.Range("E20").Value = "2026-02-16 23:19:15"
$ws.Range("N20").Value = "-4.6 °C 22:48 TU"
$ws.Range("O20").Value = "-0.9 °C"
$ws.Range("E21").Value = "2026-02-16 23:19:17"
$ws.Range("H21").Value = "66%"
$ws.Range("J21").Value = "1014.1 hPa"
$ws.Range("E22").Value = "2026-02-16 23:19:20"
$ws.Range("E23").Value = "2026-02-16 23:19:22"
$ws.Range("L23").Value = "90.7 km/h - 259º 22:58 TU"
$ws.Range("O23").Value = "-1.1 °C"
$ws.Range("E24").Value = "2026-02-16 23:19:25"
$ws.Range("O24").Value = "12.8 °C"
$ws.Range("E25").Value = "2026-02-16 23:19:27"
$ws.Range("N25").Value = "-2.6 °C 22:43 TU"
$ws.Range("O25").Value = "0.4 °C"
$ws.Range("E26").Value = "2026-02-16 23:19:30"
$ws.Range("E27").Value = "2026-02-16 23:19:32"
$ws.Range("H27").Value = "81%"
$ws.Range("N27").Value = "-1.1 °C 22:59 TU"
$ws.Range("O27").Value = "1.0 °C"
$ws.Range("E28").Value = "2026-02-16 23:19:35"
$ws.Range("J28").Value = "1012.9 hPa"
$ws.Range("O28").Value = "9.4 °C"
$ws.Range("E29").Value = "2026-02-16 23:19:37"
$ws.Range("E30").Value = "2026-02-16 23:19:39"
$ws.Range("E31").Value = "2026-02-16 23:19:42"
$ws.Range("O31").Value = "14.2 °C"
$ws.Range("E32").Value = "2026-02-16 23:19:44"
$ws.Range("E33").Value = "2026-02-16 23:19:47"
$ws.Range("E34").Value = "2026-02-16 23:19:49"
$ws.Range("O34").Value = "3.2 °C"
$ws.Range("E35").Value = "2026-02-16 23:19:52"
$ws.Range("I35").Value = "3.6 mm"
$ws.Range("J35").Value = "1016.5 hPa"
$ws.Range("E36").Value = "2026-02-16 23:19:54"
$ws.Range("H36").Value = "69%"
$ws.Range("O36").Value = "12.4 °C"
$ws.Range("E37").Value = "2026-02-16 23:19:57"
$ws.Range("O37").Value = "6.6 °C"
$ws.Range("E38").Value = "2026-02-16 23:19:59"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-16 23:20:02"
$ws.Range("L39").Value = "72.7 km/h - 352º 22:33 TU"
$ws.Range("O39").Value = "-0.2 °C"
$ws.Range("E40").Value = "2026-02-16 23:20:04"
$ws.Range("E41").Value = "2026-02-16 23:20:06"
$ws.Range("E42").Value = "2026-02-16 23:20:09"
$ws.Range("H42").Value = "78%"
$ws.Range("E43").Value = "2026-02-16 23:20:11"
$ws.Range("E44").Value = "2026-02-16 23:20:13"
$ws.Range("L44").Value = "85.7 km/h - 256º 22:34 TU"
$ws.Range("N44").Value = "-4.7 °C 22:42 TU"
$ws.Range("O44").Value = "-0.6 °C"
$ws.Range("E45").Value = "2026-02-16 23:20:16"
$ws.Range("H45").Value = "94%"
$ws.Range("E46").Value = "2026-02-16 23:20:18"
